$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 ---
$ws.Range("D2").Value = 0.8351364247873774
$ws.Range("E2").Value = 0.8351364247873774

# --- Row 3 ---
$ws.Range("D3").Value = 0.1401036061692985
$ws.Range("E3").Value = 0.1401036061692985

# --- Row 4 ---
$ws.Range("D4").Value = 0.003345906839976539
$ws.Range("E4").Value = 0.003345906839976539

# --- Row 5 ---
$ws.Range("D5").Value = 0.0001892275452282135
$ws.Range("E5").Value = 0.0001892275452282135

# --- Row 6 ---
$ws.Range("D6").Value = 0.9463777878027243
$ws.Range("E6").Value = 0.9463777878027243

# --- Row 7 ---
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0

# --- Row 8 ---
$ws.Range("D8").Value = 0.9807905980326073
$ws.Range("E8").Value = 0.01920940196739274

# --- Row 9 ---
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.01721045663602873
$ws.Range("E9").Value = 0.9827895433639713

# --- Row 10 ---
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"4.746033878922028E-06"
$ws.Range("E10").Value = 0.9999952539661211

# --- Row 11 ---
$ws.Range("D11").Value = 0.9999982070249882
$ws.Range("E11").Value = [double]"1.792975011793985E-06"
$ws.Range("F11").Value = 2.122275114059448
$ws.Range("G11").Value = 0.6

# --- Row 12 ---
$ws.Range("D12").Value = 0.897762447343373
$ws.Range("E12").Value = 0.897762447343373

# --- Row 13 ---
$ws.Range("D13").Value = 0.02225345792213777
$ws.Range("E13").Value = 0.02225345792213777

# --- Row 14 ---
$ws.Range("D14").Value = 0.001192291256906389
$ws.Range("E14").Value = 0.001192291256906389

# --- Row 15 ---
$ws.Range("D15").Value = [double]"2.617968318104033E-05"
$ws.Range("E15").Value = [double]"2.617968318104033E-05"

# --- Row 16 ---
$ws.Range("D16").Value = 0.993240162009284
$ws.Range("E16").Value = 0.993240162009284

# --- Row 17 ---
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0

# --- Row 18 ---
$ws.Range("D18").Value = 0.9999714064484416
$ws.Range("E18").Value = [double]"2.859355155837573E-05"

# --- Row 19 ---
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.00116470017827254
$ws.Range("E19").Value = 0.9988352998217275

# --- Row 20 ---
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"1.240139615713181E-10"
$ws.Range("E20").Value = 0.9999999998759861

# --- Row 21 ---
$ws.Range("D21").Value = 0.9999999991889257
$ws.Range("E21").Value = [double]"8.110743188183278E-10"
$ws.Range("F21").Value = 3.686688661575317
$ws.Range("G21").Value = 0.6

$wb.Save()
